# Automatic update of files.
# Rows 3-9 have their record data permuted/reshuffled (as re-exported by the
# source system). Destination row <- source row:
#   3 <- 5, 4 <- 6, 5 <- 8, 6 <- 4, 7 <- 3, 8 <- 9, 9 <- 7
#
# Only the columns whose values actually differ between source and
# destination rows are touched (A, B, D, E, F, G, H, P, Q, R, S); this avoids
# needlessly rewriting cells (e.g. the date/time columns) that would
# otherwise risk being reinterpreted/retyped by Excel on write-back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 3
$lastRow  = 9
$numRows  = $lastRow - $firstRow + 1   # 7

# Worksheet columns (1-based) that need to move with the permuted rows.
$cols = @(1, 2, 4, 5, 6, 7, 8, 16, 17, 18, 19)   # A, B, D, E, F, G, H, P, Q, R, S

# Mapping of destination row -> source row (worksheet row numbers).
$destToSource = @{
    3 = 5
    4 = 6
    5 = 8
    6 = 4
    7 = 3
    8 = 9
    9 = 7
}

foreach ($col in $cols) {
    $colRange = $ws.Range($ws.Cells.Item($firstRow, $col), $ws.Cells.Item($lastRow, $col))
    $colData  = $colRange.Value()   # 1-based COM array: [1..numRows, 1..1]

    $newColData = New-Object 'object[,]' $numRows, 1   # 0-based: [0..numRows-1, 0..0]

    $destRows = @(3, 4, 5, 6, 7, 8, 9)
    foreach ($destRow in $destRows) {
        $sourceRow     = $destToSource[$destRow]
        $destOffset0   = $destRow   - $firstRow        # 0-based row index into $newColData
        $sourceOffset1 = $sourceRow - $firstRow + 1    # 1-based row index into $colData
        $newColData[$destOffset0, 0] = $colData[$sourceOffset1, 1]
    }

    $colRange.Value = $newColData
}
